$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds text-formatted dates (e.g. "01-07-2021"). Setting a
# date-shaped string directly would make Excel auto-convert it to a date
# serial number, so force the cell to text format first, then clear the
# number-format override afterwards to match the existing (unstyled) data
# cells in the rest of the column.
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "01-08-2021"
$ws.Range("A45").ClearFormats()

$ws.Range("B45").Value = 110.15
$ws.Range("C45").Value = 108.41
$ws.Range("D45").Value = 111.67
$ws.Range("E45").Value = 108.36
$ws.Range("F45").Value = 117.88
